$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (GA1.9): question hint text expanded ---
$ws.Range("B10").Value = "Sort this JSON array of objects by the value of the"

# --- Row 47 (GA4.9): question hint text replaced ---
$ws.Range("B47").Value = "marks of students who scored"

# --- Row 29 (GA2.10): answer replaced by the actual ngrok URL and turned into a hyperlink ---
$ngrokUrl = "https://daeb-2409-4072-6e45-1953-c9d6-9624-b787-cecb.ngrok-free.app/"
$ws.Hyperlinks.Add($ws.Range("C29"), $ngrokUrl)
$ws.Range("C29").Value = $ngrokUrl
# Match the existing hyperlink-cell formatting used elsewhere in column C
# (copy format from another hyperlinked answer cell instead of letting a
# brand new style get appended to the stylesheet).
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# --- Column C: widen to fit the new (much longer) content ---
$ws.Columns.Item(3).AutoFit()

# --- Restore view/selection state ---
$ws.Range("B53").Select()

Write-Host "done"
